# Fabrikam Q1 marketing campaigns - header relabel + bold header row.
#
# The header row (A1:H1) carries its white font color as rich-text run
# properties inside the shared-string table; none of these cells have a
# cell-level style. A few headers also get new Korean wording. To add
# Bold while preserving that existing white run formatting, each cell's
# text is split into two Characters() sub-ranges (1..n-1 and n..n) rather
# than formatted in a single call spanning the whole string - touching
# the full string length in one Characters() call collapses the shared
# string back down to plain text under a brand-new cell-level style,
# which throws away the existing rich-text color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BoldHeader([string]$cellRef, [string]$newText) {
    $cell = $ws.Range($cellRef)

    if ($newText) {
        $cell.Value = $newText
    }

    $len = $cell.Text.Length
    $splitAt = $len - 1
    if ($splitAt -lt 1) { $splitAt = 1 }

    $first = $cell.Characters(1, $splitAt)
    $first.Font.Bold = $true
    if ($newText) { $first.Font.Color = 16777215 }

    if ($splitAt -lt $len) {
        $second = $cell.Characters($splitAt + 1, $len - $splitAt)
        $second.Font.Bold = $true
        if ($newText) { $second.Font.Color = 16777215 }
    }
}

# A1: "캠페인 담당자" -> "캠페인 소유자" (+ bold)
Set-BoldHeader "A1" "캠페인 소유자"

# C1: "시작 날짜" unchanged, just add bold
Set-BoldHeader "C1" $null

# D1: "캠페인 유형" -> "캠페인 종류" (+ bold)
Set-BoldHeader "D1" "캠페인 종류"

# G1: "총 대상 사용자" -> "총 대상 사용자 수" (+ bold)
Set-BoldHeader "G1" "총 대상 사용자 수"

# H1: "참여한 사용자" -> "참여 사용자" (+ bold)
Set-BoldHeader "H1" "참여 사용자"

# Keep the table's column names in sync with the updated header cell text
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Item(1).Name = "캠페인 소유자"
$table.ListColumns.Item(4).Name = "캠페인 종류"
$table.ListColumns.Item(7).Name = "총 대상 사용자 수"
$table.ListColumns.Item(8).Name = "참여 사용자"
